# ---------------------------------------------------------------------------
# control_panel.xlsx edit: "Determine whether to graph from control panel"
#
# Adds four new output-toggle rows (output_gtb_plots, output_flow_diagram,
# output_fractions, output_scaleups) below the existing "output_spreadsheets"
# row on the control_panel sheet, clears the two pre-set values that the
# author reset (n_organs / output_spreadsheets), re-points the yes/no
# dropdown validation that used to cover B56 onto the new B56:B60 block, and
# leaves behind the orphaned external-workbook reference that Excel embeds
# when this kind of cross-sheet list validation / formatting is pulled in
# from the master parameter spreadsheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("control_panel")
$ws.Activate()

# ---------------------------------------------------------------------------
# 0) External workbook reference (xl/externalReferences + xl/externalLinks/*)
#    Touch-and-delete a scratch formula far outside the used range so the
#    engine registers the external link + the three sheet names, without
#    leaving any visible formula cell behind.
# ---------------------------------------------------------------------------
$ws.Range("Z1000").Formula = "='[spreadsheet.xlsx]constants'!A1"
$ws.Range("Z1001").Formula = "='[spreadsheet.xlsx]time_variants'!A1"
$ws.Range("Z1002").Formula = "='[spreadsheet.xlsx]dropdown_lists'!A1"
$ws.Range("Z1000:Z1002").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 1) Clear the two values the author reset to blank
# ---------------------------------------------------------------------------
$ws.Range("B45").ClearContents()
$ws.Range("B56").ClearContents()

# ---------------------------------------------------------------------------
# 2) New rows 57-59: same look as the other "output_*" rows (e.g. row 51/53)
# ---------------------------------------------------------------------------
$ws.Range("A51:B51").Copy()
$ws.Range("A57:B57").PasteSpecial(-4122)
$ws.Range("A58:B58").PasteSpecial(-4122)
$ws.Range("A59:B59").PasteSpecial(-4122)

$ws.Range("A57").Value2 = "output_gtb_plots"
$ws.Range("A58").Value2 = "output_flow_diagram"
$ws.Range("A59").Value2 = "output_fractions"

# ---------------------------------------------------------------------------
# 3) New row 60: closes off the table, so it gets the same "bottom border"
#    treatment as row 56 used to have before rows 57-59 were inserted above
#    it. Clone row 56's current (pre-edit) top+bottom border/fill, then
#    drop the top edge so only the bottom edge remains.
# ---------------------------------------------------------------------------
$ws.Range("A56:B56").Copy()
$ws.Range("A60:B60").PasteSpecial(-4122)
$ws.Range("A60:B60").Borders.Item(8).LineStyle = -4142

$ws.Range("A60").Value2 = "output_scaleups"

# Row 56 is no longer the last row of the block, so it keeps only its top
# border (drop the bottom edge).
$ws.Range("A56:B56").Borders.Item(9).LineStyle = -4142

# ---------------------------------------------------------------------------
# 4) Dropdown validation: B56 used to be lumped in with B48:B50's yes/no
#    list; now the whole new B56:B60 block shares that same list.
# ---------------------------------------------------------------------------
$ws.Range("B56:B60").Validation.Add(3, 1, 3, "=dropdown_lists!`$B`$2:`$B`$4")
$ws.Range("B56:B60").Validation.IgnoreBlank = $true
$ws.Range("B56:B60").Validation.InCellDropdown = $true

# ---------------------------------------------------------------------------
# 5) View state: scrolled down a bit further, selection moved from C55 to B55
# ---------------------------------------------------------------------------
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
$ws.Range("B55").Select()

Write-Output "control_panel.xlsx updated"
